# Fruta / hortaliza, semanal
# A new weekly record is inserted as the new first data row (row 2),
# pushing every existing data row down by one. The former last row
# (previously row 61) becomes row 62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (first data row); this shifts
# all existing data rows (old row 2 .. old row 61) down by one, so the
# former row 61 becomes row 62 automatically.
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits formatting from the row below; clear
# it so the new row matches the plain (unstyled) look of the other data
# rows before we reapply the one style that is actually used (date format
# on column D).
$ws.Range("A2:R2").ClearFormats()

# Populate the new record (week of 2022-03-23).
$ws.Range("A2").Value() = 1
$ws.Range("B2").Value() = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value() = "Arica y Parinacota"
$ws.Range("D2").Value() = 44643
$ws.Range("E2").Value() = 15
$ws.Range("F2").Value() = 100112012
$ws.Range("G2").Value() = "Espinaca"
$ws.Range("H2").Value() = "Sin especificar"
$ws.Range("I2").Value() = "Primera"
$ws.Range("J2").Value() = 300
$ws.Range("K2").Value() = 900
$ws.Range("L2").Value() = 1000
$ws.Range("M2").Value() = 950
$ws.Range("N2").Value() = "$/atado 2,5 a 3 kilos"
$ws.Range("O2").Value() = "Región de Arica y Parinacota"
$ws.Range("P2").Value() = 317
$ws.Range("Q2").Value() = 3
$ws.Range("R2").Value() = "Hortaliza"

# Column D holds dates; restore the custom date number format used
# throughout the rest of the sheet.
$ws.Range("D2").NumberFormat() = "YYYY-MM-DD HH:MM:SS"
